$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update item rows (8-22): quantities, serial numbers, descriptions, rates and amounts
# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = ''
$ws.Range("C8").Value = 37
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.0'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '0.00'

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = ''
$ws.Range("C9").Value = 78
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.0'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F9").Value = 0
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '0.00'

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'P. point'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = 'On board'
$ws.Range("F10").Value = 136
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '12376.00'

# Row 11
$ws.Range("C11").Value = 71
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.0'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 23
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '1633.00'

# Row 12
$ws.Range("C12").Value = 22
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.0'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 50
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '1100.00'

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = 'Each'
$ws.Range("C13").Value = 32
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.0'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 219
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '7008.00'

# Row 14
$ws.Range("C14").Value = 82
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.0'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = 'R. mtr.'
$ws.Range("C15").Value = 87
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '17'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '25 mm'
$ws.Range("F15").Value = 56
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '4872.00'

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = 'Mtr.'
$ws.Range("C16").Value = 25
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F16").Value = 122
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '3050.00'

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = 'Mtr.'
$ws.Range("C17").Value = 22
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F17").Value = 20
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '440.00'

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = ''
$ws.Range("C18").Value = 69
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.0'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F18").Value = 0
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '0.00'

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = 'Each'
$ws.Range("C19").Value = 2
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '30'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = ' 6 A to 32 A rating'
$ws.Range("F19").Value = 187
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '374.00'

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = ''
$ws.Range("C20").Value = 63
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.0'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F20").Value = 0
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '0.00'

# Row 21
$ws.Range("C21").Value = 70

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = ''
$ws.Range("C22").Value = 81
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '38'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = 'Grand Total'

# Remove the obsolete "Grand Total" row (old row 23); rows below shift up by one
$ws.Range("A23").EntireRow.Delete()

# Update recalculated totals (now on rows 24 and 26 after the shift)
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "30853.00"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "30853.00"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "30853.00"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "30853.00"
